$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New benchmark run was added to the source data: one more column per
#     group (G/N/U), holding the new "reparse removed" timings. -----------
$ws.Range("G2").Value = 5.7359999999999998
$ws.Range("G3").Value = 4.8719999999999999

$ws.Range("N2").Value = 11.29
$ws.Range("N3").Value = 10.631

$ws.Range("U2").Value = 17.010000000000002
$ws.Range("U3").Value = 16.853999999999999

# --- Extend the three line-chart series so they plot the new column too --
$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$A`$2,,Sheet1!`$B`$2:`$G`$2,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES(Sheet1!`$A`$3,,Sheet1!`$B`$3:`$G`$3,2)"

$chart2 = $ws.ChartObjects(2).Chart
$chart2.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$H`$2,,Sheet1!`$I`$2:`$N`$2,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES(Sheet1!`$H`$3,,Sheet1!`$I`$3:`$N`$3,2)"

$chart3 = $ws.ChartObjects(3).Chart
$chart3.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$O`$2,,Sheet1!`$P`$2:`$U`$2,1)"
$chart3.SeriesCollection(2).Formula = "=SERIES(Sheet1!`$O`$3,,Sheet1!`$P`$3:`$U`$3,2)"

# --- The three chart frames were each trimmed a little shorter (and the
#     third, narrower too) once the extra series point was plotted. ------
$co1 = $ws.ChartObjects(1)
$co1.Height = 345.74992125984255

$co2 = $ws.ChartObjects(2)
$co2.Height = 343.5000393700785

$co3 = $ws.ChartObjects(3)
$co3.Height = 344.25011811023654
$co3.Width = 423.42281219242153

# --- Selection cursor ended up on Y18 -------------------------------------
$ws.Range("Y18").Select()
